$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.488.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.663.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.80%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.123"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.413"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.87%  "

$ws.Range("E12").Value = "  +1.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.38%  "

$ws.Range("E14").Value = "  +20.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.137.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.251.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.661.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "360.89"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.72%  "

$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.57%  "

$ws.Range("E24").Value = "  +1.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000102"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +16.96%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.45%  "

$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "558.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.73%  "

$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.166"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.44%  "

$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.435"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "163.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("E42").Value = "  +7.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "168.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0627"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.665"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0265"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0990"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.40%  "

